$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.711.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.09%  "
$ws.Range("D3").Value = "'3.157.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.24%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'525.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'133.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.13%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'3.156.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.13%  "
$ws.Range("E9").Value = "  -6.32%  "
$ws.Range("E10").Value = "  -6.86%  "
$ws.Range("E11").Value = "  -6.61%  "
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("D13").Value = "'3.701.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.33%  "
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "'25.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.84%  "
$ws.Range("D16").Value = "'3.160.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.72%  "
$ws.Range("D17").Value = "'57.723.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.12%  "
$ws.Range("E18").Value = "  -7.59%  "
$ws.Range("E19").Value = "  -4.63%  "
$ws.Range("E20").Value = "  -8.66%  "
$ws.Range("D21").Value = "'8.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.56%  "
$ws.Range("D22").Value = "'346.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.17%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'69.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.74%  "
$ws.Range("E25").Value = "  -6.56%  "
$ws.Range("D26").Value = "'3.289.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("D27").Value = "'0.0₃0960"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.14%  "
$ws.Range("E28").Value = "  -3.98%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'6.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.36%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("E32").Value = "  -7.79%  "
$ws.Range("D33").Value = "'6.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.19%  "
$ws.Range("D34").Value = "'21.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.57%  "
$ws.Range("D35").Value = "'1.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.46%  "
$ws.Range("D36").Value = "'4.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.07%  "
$ws.Range("D37").Value = "'159.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("E38").Value = "  -7.21%  "
$ws.Range("E39").Value = "  -7.76%  "
$ws.Range("D40").Value = "'25.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.98%  "
$ws.Range("D41").Value = "'0.0696"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.18%  "
$ws.Range("D42").Value = "'3.185.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("D43").Value = "'40.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("E44").Value = "  -6.65%  "
$ws.Range("E45").Value = "  -3.49%  "
$ws.Range("E46").Value = "  -6.02%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'1.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.17%  "
$ws.Range("D49").Value = "'2.269.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("D50").Value = "'6.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.90%  "
$ws.Range("D51").Value = "'20.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.60%  "
